$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("alpha_school_climate")
$ws.Range("B2").Value = 0.93476056157324883
$ws.Range("B5").Value = 0.52373588132351678
$ws.Range("C5").Value = 0.46357909123400726
$ws.Range("D5").Value = 0.42846196583932128
$ws.Range("E5").Value = 0.93439883556814018
$ws.Range("B6").Value = 0.59251316719084646
$ws.Range("C6").Value = 0.53891351314369262
$ws.Range("D6").Value = 0.42333952995152985
$ws.Range("E6").Value = 0.93310289925705125
$ws.Range("B7").Value = 0.47399432214396608
$ws.Range("C7").Value = 0.41091850756826753
$ws.Range("D7").Value = 0.43250728611586581
$ws.Range("E7").Value = 0.93540304587500611
$ws.Range("B8").Value = 0.54909354335020133
$ws.Range("C8").Value = 0.49152066719861881
$ws.Range("D8").Value = 0.4267227188967791
$ws.Range("E8").Value = 0.93396190638022281
$ws.Range("B9").Value = 0.51890870989776039
$ws.Range("C9").Value = 0.45715157956819463
$ws.Range("D9").Value = 0.42996012074238438
$ws.Range("E9").Value = 0.93477268865455765
$ws.Range("B10").Value = 0.77090251752710426
$ws.Range("C10").Value = 0.73492784796474975
$ws.Range("D10").Value = 0.40927159241640076
$ws.Range("E10").Value = 0.92939688439252799
$ws.Range("B11").Value = 0.80061650505030113
$ws.Range("C11").Value = 0.76826267875016951
$ws.Range("D11").Value = 0.40673690572357024
$ws.Range("E11").Value = 0.92870517355628091
$ws.Range("B12").Value = 0.75210045074974896
$ws.Range("C12").Value = 0.71299111933680426
$ws.Range("D12").Value = 0.41044119372959637
$ws.Range("E12").Value = 0.92971352897790249
$ws.Range("B13").Value = 0.83990682197243183
$ws.Range("C13").Value = 0.81318604897136504
$ws.Range("D13").Value = 0.40355781856995915
$ws.Range("E13").Value = 0.92782681007882262
$ws.Range("B14").Value = 0.73539096330184228
$ws.Range("C14").Value = 0.69614789143724476
$ws.Range("D14").Value = 0.41199258273379419
$ws.Range("E14").Value = 0.93013109012361239
$ws.Range("B15").Value = 0.76595155849930652
$ws.Range("C15").Value = 0.72982799609508942
$ws.Range("D15").Value = 0.40952865767137958
$ws.Range("E15").Value = 0.92946661587351143
$ws.Range("B16").Value = 0.81908087946792407
$ws.Range("C16").Value = 0.79026917878875036
$ws.Range("D16").Value = 0.40533097345384811
$ws.Range("E16").Value = 0.92831821704040562
$ws.Range("B17").Value = 0.79621466286102838
$ws.Range("C17").Value = 0.76417205183184378
$ws.Range("D17").Value = 0.40721760749368463
$ws.Range("E17").Value = 0.92883693833443903
$ws.Range("B18").Value = 0.27702216209658148
$ws.Range("C18").Value = 0.20053095095819762
$ws.Range("D18").Value = 0.44906847625775009
$ws.Range("E18").Value = 0.93934631538536573
$ws.Range("B19").Value = 0.65990316316536268
$ws.Range("C19").Value = 0.61395412342376854
$ws.Range("D19").Value = 0.41740737239365239
$ws.Range("E19").Value = 0.93156703162849719
$ws.Range("B20").Value = 0.60744512371337112
$ws.Range("C20").Value = 0.55265868503004423
$ws.Range("D20").Value = 0.42206164892062714
$ws.Range("E20").Value = 0.93277527302164553
$ws.Range("B21").Value = 0.74400183617474536
$ws.Range("C21").Value = 0.70473948744544235
$ws.Range("D21").Value = 0.41125401993791555
$ws.Range("E21").Value = 0.9299326497986401
$ws.Range("B22").Value = 0.75865434474670312
$ws.Range("C22").Value = 0.72375332757437916
$ws.Range("D22").Value = 0.41098419044898116
$ws.Range("E22").Value = 0.92985999411442655
$ws.Range("B23").Value = 0.6762775285988154
$ws.Range("C23").Value = 0.63007304074528725
$ws.Range("D23").Value = 0.41651398454986593
$ws.Range("E23").Value = 0.93133238433475152
$ws.Range("B24").Value = 0.69386827898254022
$ws.Range("C24").Value = 0.64837663663662637
$ws.Range("D24").Value = 0.41536514361429455
$ws.Range("E24").Value = 0.93102933402395294

$ws = $wb.Worksheets.Item("alpha_teacher_quality")
$ws.Range("B2").Value = 0.89085471361440993
$ws.Range("B5").Value = 0.60819780597281548
$ws.Range("C5").Value = 0.53935298750473015
$ws.Range("D5").Value = 0.3238364661916715
$ws.Range("E5").Value = 0.88456541029362779
$ws.Range("B6").Value = 0.57778376077879723
$ws.Range("C6").Value = 0.50276421382225001
$ws.Range("D6").Value = 0.32606512477130556
$ws.Range("E6").Value = 0.88559879118161744
$ws.Range("B7").Value = 0.56995610225589366
$ws.Range("C7").Value = 0.49766811609716827
$ws.Range("D7").Value = 0.32638993295791435
$ws.Range("E7").Value = 0.88574841932256598
$ws.Range("B8").Value = 0.53645401598006992
$ws.Range("C8").Value = 0.45654114235215276
$ws.Range("D8").Value = 0.32961501160310597
$ws.Range("E8").Value = 0.88722079518052255
$ws.Range("B9").Value = 0.54773926891713187
$ws.Range("C9").Value = 0.47202409141936263
$ws.Range("D9").Value = 0.32849871559948751
$ws.Range("E9").Value = 0.88671388157372577
$ws.Range("B10").Value = 0.72941299464206544
$ws.Range("C10").Value = 0.67393034620805248
$ws.Range("D10").Value = 0.31310874268859251
$ws.Range("E10").Value = 0.87942149660198099
$ws.Range("B11").Value = 0.7555214115293436
$ws.Range("C11").Value = 0.70509800240443421
$ws.Range("D11").Value = 0.31079282300349548
$ws.Range("E11").Value = 0.87827264543474437
$ws.Range("B12").Value = 0.71592779032245291
$ws.Range("C12").Value = 0.659755437365384
$ws.Range("D12").Value = 0.3142985261712421
$ws.Range("E12").Value = 0.88000627761278505
$ws.Range("B13").Value = 0.68707030680333925
$ws.Range("C13").Value = 0.62678762214784201
$ws.Range("D13").Value = 0.31695771008770446
$ws.Range("E13").Value = 0.88130015864456968
$ws.Range("B14").Value = 0.74111284797246979
$ws.Range("C14").Value = 0.68819618466642296
$ws.Range("D14").Value = 0.31204809528848637
$ws.Range("E14").Value = 0.8788970884648748
$ws.Range("B15").Value = 0.68861765846443657
$ws.Range("C15").Value = 0.62798229627429736
$ws.Range("D15").Value = 0.31689179019138158
$ws.Range("E15").Value = 0.88126830081228347
$ws.Range("B16").Value = 0.44039033537765104
$ws.Range("C16").Value = 0.35068942179869245
$ws.Range("D16").Value = 0.33943492711405776
$ws.Range("E16").Value = 0.89155996381164715
$ws.Range("B17").Value = 0.46516125289689464
$ws.Range("C17").Value = 0.37744053680909989
$ws.Range("D17").Value = 0.33712299320784628
$ws.Range("E17").Value = 0.89055737468229412
$ws.Range("B18").Value = 0.5547657433496157
$ws.Range("C18").Value = 0.47615009455247082
$ws.Range("D18").Value = 0.32930500997207696
$ws.Range("E18").Value = 0.88708030908991553
$ws.Range("B19").Value = 0.50651347460288754
$ws.Range("C19").Value = 0.42256245024194994
$ws.Range("D19").Value = 0.33320641973330223
$ws.Range("E19").Value = 0.88883246494201507
$ws.Range("B20").Value = 0.49238339463829822
$ws.Range("C20").Value = 0.40741265810670729
$ws.Range("D20").Value = 0.33493932461821685
$ws.Range("E20").Value = 0.8895998084666662
$ws.Range("B21").Value = 0.63608330714518779
$ws.Range("C21").Value = 0.56675496613262022
$ws.Range("D21").Value = 0.32179368951610843
$ws.Range("E21").Value = 0.88360780369828562

$ws = $wb.Worksheets.Item("alpha_student_support")
$ws.Range("B2").Value = 0.65623927279165772
$ws.Range("B5").Value = 0.72343446334671924
$ws.Range("C5").Value = 0.45186009875595773
$ws.Range("D5").Value = 0.30679225020535938
$ws.Range("E5").Value = 0.57039265708229547
$ws.Range("B6").Value = 0.6483559048496802
$ws.Range("C6").Value = 0.34956194140277003
$ws.Range("D6").Value = 0.37960635378077401
$ws.Range("E6").Value = 0.64734585900123698
$ws.Range("B7").Value = 0.74150754987620371
$ws.Range("C7").Value = 0.45718989533615068
$ws.Range("D7").Value = 0.298158036314709
$ws.Range("E7").Value = 0.56033646737063059
$ws.Range("B8").Value = 0.72253835183014314
$ws.Range("C8").Value = 0.45330117956787463
$ws.Range("D8").Value = 0.30698763860471773
$ws.Range("E8").Value = 0.57061773424838258

$ws = $wb.Worksheets.Item("alpha_student_motivation")
$ws.Range("B2").Value = 0.9534737653335531
$ws.Range("B5").Value = 0.9400258716737856
$ws.Range("C5").Value = 0.89202801180388513
$ws.Range("D5").Value = 0.83224411654333985
$ws.Range("E5").Value = 0.93704011097758799
$ws.Range("B6").Value = 0.88296490018819718
$ws.Range("C6").Value = 0.79537574322487448
$ws.Range("D6").Value = 0.90487469804358511
$ws.Range("E6").Value = 0.96614456004906146
$ws.Range("B7").Value = 0.96712286367211653
$ws.Range("C7").Value = 0.94006630972687655
$ws.Range("D7").Value = 0.79860742705905385
$ws.Range("E7").Value = 0.92245825461007924
$ws.Range("B8").Value = 0.95834418544591693
$ws.Range("C8").Value = 0.92380541615946565
$ws.Range("D8").Value = 0.81104528041326462
$ws.Range("E8").Value = 0.9279373785140459
